$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "('MRV_Method', 'forward_method')"
$ws.Range("C2").Value = 3.510622613430023

# Row 3
$ws.Range("B3").Value = "('MRV_Degree_Method', 'forward_method')"
$ws.Range("C3").Value = 4.016371684074402

# Row 4
$ws.Range("B4").Value = "('MRV_Degree_Method', 'ac3_method')"
$ws.Range("C4").Value = 4.982024788856506

# Row 5
$ws.Range("B5").Value = "('MRV_Method', 'ac3_method')"
$ws.Range("C5").Value = 4.292513732910156
